$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1138-1139; this shifts the existing rows
# 1138-1230 down to 1140-1232 and extends the used range to A1:R1232.
$ws.Rows("1138:1139").Insert()

# Populate the two newly-inserted rows with the new weekly data.
# Row 1138
$ws.Cells.Item(1138, 1).Value = 6
$ws.Cells.Item(1138, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1138, 3).Value = "Metropolitana"
$ws.Cells.Item(1138, 4).Value = 45013
$ws.Cells.Item(1138, 5).Value = 13
$ws.Cells.Item(1138, 6).Value = 100112040
$ws.Cells.Item(1138, 7).Value = "Cilantro"
$ws.Cells.Item(1138, 8).Value = "Sin especificar"
$ws.Cells.Item(1138, 9).Value = "Primera"
$ws.Cells.Item(1138, 10).Value = 680
$ws.Cells.Item(1138, 11).Value = 5000
$ws.Cells.Item(1138, 12).Value = 6000
$ws.Cells.Item(1138, 13).Value = 5515
$ws.Cells.Item(1138, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(1138, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1138, 16).Value = 153
$ws.Cells.Item(1138, 17).Value = 36
$ws.Cells.Item(1138, 18).Value = "Hortaliza"

# Row 1139
$ws.Cells.Item(1139, 1).Value = 6
$ws.Cells.Item(1139, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1139, 3).Value = "Metropolitana"
$ws.Cells.Item(1139, 4).Value = 45013
$ws.Cells.Item(1139, 5).Value = 13
$ws.Cells.Item(1139, 6).Value = 100112040
$ws.Cells.Item(1139, 7).Value = "Cilantro"
$ws.Cells.Item(1139, 8).Value = "Sin especificar"
$ws.Cells.Item(1139, 9).Value = "Primera"
$ws.Cells.Item(1139, 10).Value = 510
$ws.Cells.Item(1139, 11).Value = 10000
$ws.Cells.Item(1139, 12).Value = 11000
$ws.Cells.Item(1139, 13).Value = 10431
$ws.Cells.Item(1139, 14).Value = "$/docena de atados"
$ws.Cells.Item(1139, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1139, 16).Value = 3477
$ws.Cells.Item(1139, 17).Value = 3
$ws.Cells.Item(1139, 18).Value = "Hortaliza"
